$wb = $excel.ActiveWorkbook

# --- Sheet "展览" (Exhibition) ---
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F3").Value = 599
$ws1.Range("F5").Value = 585
$ws1.Range("F9").Value = 203
$ws1.Range("F11").Value = 484
$ws1.Range("F12").Value = 1438
$ws1.Range("F14").Value = 132
$ws1.Range("F17").Value = 111
$ws1.Range("F18").Value = 691
$ws1.Range("F19").Value = 1030
$ws1.Range("F20").Value = 48
$ws1.Range("F21").Value = 282
$ws1.Range("F23").Value = 6186
$ws1.Range("F24").Value = 76
$ws1.Range("F25").Value = 133
$ws1.Range("F26").Value = 131
$ws1.Range("F28").Value = 14959
$ws1.Range("F30").Value = 257
$ws1.Range("F33").Value = 10893
$ws1.Range("F34").Value = 699
$ws1.Range("F35").Value = 4261
$ws1.Range("G7").Value = "不可售"

# --- Sheet "全部类型" (All types) ---
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F3").Value = 599
$ws4.Range("F5").Value = 585
$ws4.Range("F9").Value = 203
$ws4.Range("F11").Value = 484
$ws4.Range("F12").Value = 1438
$ws4.Range("F14").Value = 132
$ws4.Range("F18").Value = 111
$ws4.Range("F19").Value = 691
$ws4.Range("F21").Value = 1030
$ws4.Range("F22").Value = 48
$ws4.Range("F23").Value = 282
$ws4.Range("F26").Value = 6186
$ws4.Range("F27").Value = 76
$ws4.Range("F28").Value = 133
$ws4.Range("F29").Value = 131
$ws4.Range("F31").Value = 14959
$ws4.Range("F33").Value = 257
$ws4.Range("F36").Value = 10893
$ws4.Range("F37").Value = 699
$ws4.Range("F38").Value = 4261
$ws4.Range("G7").Value = "不可售"
